# Apply the blueprint_task_products.xlsx changes:
#  - Rename "roomplan FP" -> "roomplan-FP" (cell A7)
#  - Rename "Walkthrough video" -> "Walkthrough-video" (cell A23)
#  - Update the view: zoom from 100% to 160%, and active selection from A8 to A7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the task name values
$ws.Range("A7").Value = "roomplan-FP"
$ws.Range("A23").Value = "Walkthrough-video"

# Update zoom level in the active window view
$excel.ActiveWindow.Zoom = 160

# Update the active cell / selection to A7
$ws.Range("A7").Select()
